$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.9792566666666667
$ws.Cells.Item(2, 8).Value = 2.93777
$ws.Cells.Item(2, 9).Value = 0.02840115057834171
$ws.Cells.Item(2, 10).Value = 0.02840115057834171
$ws.Cells.Item(2, 13).Value = 10.08846466666667
$ws.Cells.Item(2, 14).Value = 30.265394
$ws.Cells.Item(2, 15).Value = 0.2597171077778241
$ws.Cells.Item(2, 16).Value = 0.2597171077778241
$ws.Cells.Item(2, 17).Value = 9.879196281264445
$ws.Cells.Item(2, 18).Value = 88.91276653138
$ws.Cells.Item(2, 19).Value = 0.007376264685769388
$ws.Cells.Item(2, 20).Value = 0.007376264685769387
$ws.Cells.Item(3, 7).Value = 0.9792566666666667
$ws.Cells.Item(3, 8).Value = 2.93777
$ws.Cells.Item(3, 9).Value = 0.02840115057834171
$ws.Cells.Item(3, 10).Value = 0.02840115057834171
$ws.Cells.Item(3, 15).Value = 0.05876531725312483
$ws.Cells.Item(3, 16).Value = 0.05876531725312482
$ws.Cells.Item(3, 17).Value = 2.235332545636667
$ws.Cells.Item(3, 18).Value = 20.11799291073
$ws.Cells.Item(3, 19).Value = 0.001669002624090021
$ws.Cells.Item(3, 20).Value = 0.00166900262409002
$ws.Cells.Item(4, 7).Value = 0.9792566666666667
$ws.Cells.Item(4, 8).Value = 2.93777
$ws.Cells.Item(4, 9).Value = 0.02840115057834171
$ws.Cells.Item(4, 10).Value = 0.02840115057834171
$ws.Cells.Item(4, 13).Value = 14.90894133333333
$ws.Cells.Item(4, 14).Value = 44.726824
$ws.Cells.Item(4, 15).Value = 0.3838153030278664
$ws.Cells.Item(4, 16).Value = 0.3838153030278664
$ws.Cells.Item(4, 17).Value = 14.59968019360889
$ws.Cells.Item(4, 18).Value = 131.39712174248
$ws.Cells.Item(4, 19).Value = 0.01090079621556629
$ws.Cells.Item(4, 20).Value = 0.01090079621556629
$ws.Cells.Item(5, 7).Value = 0.9792566666666667
$ws.Cells.Item(5, 8).Value = 2.93777
$ws.Cells.Item(5, 9).Value = 0.02840115057834171
$ws.Cells.Item(5, 10).Value = 0.02840115057834171
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 0.2836386666666667
$ws.Cells.Item(5, 14).Value = 0.850916
$ws.Cells.Item(5, 15).Value = 0.007301984652235982
$ws.Cells.Item(5, 16).Value = 0.007301984652235982
$ws.Cells.Item(5, 17).Value = 0.2777550552577778
$ws.Cells.Item(5, 18).Value = 2.49979549732
$ws.Cells.Item(5, 19).Value = 0.0002073847656288943
$ws.Cells.Item(5, 20).Value = 0.0002073847656288943
$ws.Cells.Item(6, 7).Value = 0.9792566666666667
$ws.Cells.Item(6, 8).Value = 2.93777
$ws.Cells.Item(6, 9).Value = 0.02840115057834171
$ws.Cells.Item(6, 10).Value = 0.02840115057834171
$ws.Cells.Item(6, 13).Value = 11.28032366666667
$ws.Cells.Item(6, 14).Value = 33.840971
$ws.Cells.Item(6, 15).Value = 0.2904002872889486
$ws.Cells.Item(6, 16).Value = 0.2904002872889486
$ws.Cells.Item(6, 17).Value = 11.04633215274111
$ws.Cells.Item(6, 18).Value = 99.41698937466998
$ws.Cells.Item(6, 19).Value = 0.008247702287287122
$ws.Cells.Item(6, 20).Value = 0.00824770228728712
$ws.Cells.Item(7, 9).Value = 0.06296007145894493
$ws.Cells.Item(7, 10).Value = 0.06296007145894492
$ws.Cells.Item(7, 13).Value = 10.08846466666667
$ws.Cells.Item(7, 14).Value = 30.265394
$ws.Cells.Item(7, 15).Value = 0.2597171077778241
$ws.Cells.Item(7, 16).Value = 0.2597171077778241
$ws.Cells.Item(7, 17).Value = 21.90034175234
$ws.Cells.Item(7, 18).Value = 197.10307577106
$ws.Cells.Item(7, 19).Value = 0.01635180766480231
$ws.Cells.Item(7, 20).Value = 0.01635180766480231
$ws.Cells.Item(8, 9).Value = 0.06296007145894493
$ws.Cells.Item(8, 10).Value = 0.06296007145894492
$ws.Cells.Item(8, 15).Value = 0.05876531725312483
$ws.Cells.Item(8, 16).Value = 0.05876531725312482
$ws.Cells.Item(8, 19).Value = 0.003699868573564309
$ws.Cells.Item(8, 20).Value = 0.003699868573564307
$ws.Cells.Item(9, 9).Value = 0.06296007145894493
$ws.Cells.Item(9, 10).Value = 0.06296007145894492
$ws.Cells.Item(9, 13).Value = 14.90894133333333
$ws.Cells.Item(9, 14).Value = 44.726824
$ws.Cells.Item(9, 15).Value = 0.3838153030278664
$ws.Cells.Item(9, 16).Value = 0.3838153030278664
$ws.Cells.Item(9, 17).Value = 32.36477711464
$ws.Cells.Item(9, 18).Value = 291.28299403176
$ws.Cells.Item(9, 19).Value = 0.02416503890567107
$ws.Cells.Item(9, 20).Value = 0.02416503890567107
$ws.Cells.Item(10, 9).Value = 0.06296007145894493
$ws.Cells.Item(10, 10).Value = 0.06296007145894492
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 0.2836386666666667
$ws.Cells.Item(10, 14).Value = 0.850916
$ws.Cells.Item(10, 15).Value = 0.007301984652235982
$ws.Cells.Item(10, 16).Value = 0.007301984652235982
$ws.Cells.Item(10, 17).Value = 0.61573132676
$ws.Cells.Item(10, 18).Value = 5.54158194084
$ws.Cells.Item(10, 19).Value = 0.0004597334754968966
$ws.Cells.Item(10, 20).Value = 0.0004597334754968965
$ws.Cells.Item(11, 9).Value = 0.06296007145894493
$ws.Cells.Item(11, 10).Value = 0.06296007145894492
$ws.Cells.Item(11, 13).Value = 11.28032366666667
$ws.Cells.Item(11, 14).Value = 33.840971
$ws.Cells.Item(11, 15).Value = 0.2904002872889486
$ws.Cells.Item(11, 16).Value = 0.2904002872889486
$ws.Cells.Item(11, 17).Value = 24.48766502531
$ws.Cells.Item(11, 18).Value = 220.38898522779
$ws.Cells.Item(11, 19).Value = 0.01828362283941034
$ws.Cells.Item(11, 20).Value = 0.01828362283941034
$ws.Cells.Item(12, 7).Value = 18.980972
$ws.Cells.Item(12, 8).Value = 56.942916
$ws.Cells.Item(12, 9).Value = 0.550500662640664
$ws.Cells.Item(12, 10).Value = 0.550500662640664
$ws.Cells.Item(12, 13).Value = 10.08846466666667
$ws.Cells.Item(12, 14).Value = 30.265394
$ws.Cells.Item(12, 15).Value = 0.2597171077778241
$ws.Cells.Item(12, 16).Value = 0.2597171077778241
$ws.Cells.Item(12, 17).Value = 191.4888653609893
$ws.Cells.Item(12, 18).Value = 1723.399788248904
$ws.Cells.Item(12, 19).Value = 0.1429744399308089
$ws.Cells.Item(12, 20).Value = 0.1429744399308089
$ws.Cells.Item(13, 7).Value = 18.980972
$ws.Cells.Item(13, 8).Value = 56.942916
$ws.Cells.Item(13, 9).Value = 0.550500662640664
$ws.Cells.Item(13, 10).Value = 0.550500662640664
$ws.Cells.Item(13, 15).Value = 0.05876531725312483
$ws.Cells.Item(13, 16).Value = 0.05876531725312482
$ws.Cells.Item(13, 17).Value = 43.32754210787599
$ws.Cells.Item(13, 18).Value = 389.9478789708839
$ws.Cells.Item(13, 19).Value = 0.03235034608813406
$ws.Cells.Item(13, 20).Value = 0.03235034608813406
$ws.Cells.Item(14, 7).Value = 18.980972
$ws.Cells.Item(14, 8).Value = 56.942916
$ws.Cells.Item(14, 9).Value = 0.550500662640664
$ws.Cells.Item(14, 10).Value = 0.550500662640664
$ws.Cells.Item(14, 13).Value = 14.90894133333333
$ws.Cells.Item(14, 14).Value = 44.726824
$ws.Cells.Item(14, 15).Value = 0.3838153030278664
$ws.Cells.Item(14, 16).Value = 0.3838153030278664
$ws.Cells.Item(14, 17).Value = 282.9861979976426
$ws.Cells.Item(14, 18).Value = 2546.875781978784
$ws.Cells.Item(14, 19).Value = 0.2112905786484677
$ws.Cells.Item(14, 20).Value = 0.2112905786484677
$ws.Cells.Item(15, 7).Value = 18.980972
$ws.Cells.Item(15, 8).Value = 56.942916
$ws.Cells.Item(15, 9).Value = 0.550500662640664
$ws.Cells.Item(15, 10).Value = 0.550500662640664
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 0.2836386666666667
$ws.Cells.Item(15, 14).Value = 0.850916
$ws.Cells.Item(15, 15).Value = 0.007301984652235982
$ws.Cells.Item(15, 16).Value = 0.007301984652235982
$ws.Cells.Item(15, 17).Value = 5.383737590117333
$ws.Cells.Item(15, 18).Value = 48.453638311056
$ws.Cells.Item(15, 19).Value = 0.004019747389647867
$ws.Cells.Item(15, 20).Value = 0.004019747389647867
$ws.Cells.Item(16, 7).Value = 18.980972
$ws.Cells.Item(16, 8).Value = 56.942916
$ws.Cells.Item(16, 9).Value = 0.550500662640664
$ws.Cells.Item(16, 10).Value = 0.550500662640664
$ws.Cells.Item(16, 13).Value = 11.28032366666667
$ws.Cells.Item(16, 14).Value = 33.840971
$ws.Cells.Item(16, 15).Value = 0.2904002872889486
$ws.Cells.Item(16, 16).Value = 0.2904002872889486
$ws.Cells.Item(16, 17).Value = 214.1115076679373
$ws.Cells.Item(16, 18).Value = 1927.003569011436
$ws.Cells.Item(16, 19).Value = 0.1598655505836054
$ws.Cells.Item(16, 20).Value = 0.1598655505836054
$ws.Cells.Item(17, 7).Value = 0.9440396666666667
$ws.Cells.Item(17, 8).Value = 2.832119
$ws.Cells.Item(17, 9).Value = 0.02737976021771022
$ws.Cells.Item(17, 10).Value = 0.02737976021771022
$ws.Cells.Item(17, 13).Value = 10.08846466666667
$ws.Cells.Item(17, 14).Value = 30.265394
$ws.Cells.Item(17, 15).Value = 0.2597171077778241
$ws.Cells.Item(17, 16).Value = 0.2597171077778241
$ws.Cells.Item(17, 17).Value = 9.523910821098445
$ws.Cells.Item(17, 18).Value = 85.715197389886
$ws.Cells.Item(17, 19).Value = 0.007110992135394028
$ws.Cells.Item(17, 20).Value = 0.007110992135394027
$ws.Cells.Item(18, 7).Value = 0.9440396666666667
$ws.Cells.Item(18, 8).Value = 2.832119
$ws.Cells.Item(18, 9).Value = 0.02737976021771022
$ws.Cells.Item(18, 10).Value = 0.02737976021771022
$ws.Cells.Item(18, 15).Value = 0.05876531725312483
$ws.Cells.Item(18, 16).Value = 0.05876531725312482
$ws.Cells.Item(18, 17).Value = 2.154943298425667
$ws.Cells.Item(18, 18).Value = 19.394489685831
$ws.Cells.Item(18, 19).Value = 0.001608980295508228
$ws.Cells.Item(18, 20).Value = 0.001608980295508227
$ws.Cells.Item(19, 7).Value = 0.9440396666666667
$ws.Cells.Item(19, 8).Value = 2.832119
$ws.Cells.Item(19, 9).Value = 0.02737976021771022
$ws.Cells.Item(19, 10).Value = 0.02737976021771022
$ws.Cells.Item(19, 13).Value = 14.90894133333333
$ws.Cells.Item(19, 14).Value = 44.726824
$ws.Cells.Item(19, 15).Value = 0.3838153030278664
$ws.Cells.Item(19, 16).Value = 0.3838153030278664
$ws.Cells.Item(19, 17).Value = 14.07463200667289
$ws.Cells.Item(19, 18).Value = 126.671688060056
$ws.Cells.Item(19, 19).Value = 0.01050877096479077
$ws.Cells.Item(19, 20).Value = 0.01050877096479077
$ws.Cells.Item(20, 7).Value = 0.9440396666666667
$ws.Cells.Item(20, 8).Value = 2.832119
$ws.Cells.Item(20, 9).Value = 0.02737976021771022
$ws.Cells.Item(20, 10).Value = 0.02737976021771022
$ws.Cells.Item(20, 11).Value = 3.0
$ws.Cells.Item(20, 12).Value = 1.0
$ws.Cells.Item(20, 13).Value = 0.2836386666666667
$ws.Cells.Item(20, 14).Value = 0.850916
$ws.Cells.Item(20, 15).Value = 0.007301984652235982
$ws.Cells.Item(20, 16).Value = 0.007301984652235982
$ws.Cells.Item(20, 17).Value = 0.2677661523337778
$ws.Cells.Item(20, 18).Value = 2.409895371004
$ws.Cells.Item(20, 19).Value = 0.0001999265888916214
$ws.Cells.Item(20, 20).Value = 0.0001999265888916213
$ws.Cells.Item(21, 7).Value = 0.9440396666666667
$ws.Cells.Item(21, 8).Value = 2.832119
$ws.Cells.Item(21, 9).Value = 0.02737976021771022
$ws.Cells.Item(21, 10).Value = 0.02737976021771022
$ws.Cells.Item(21, 13).Value = 11.28032366666667
$ws.Cells.Item(21, 14).Value = 33.840971
$ws.Cells.Item(21, 15).Value = 0.2904002872889486
$ws.Cells.Item(21, 16).Value = 0.2904002872889486
$ws.Cells.Item(21, 17).Value = 10.64907299417211
$ws.Cells.Item(21, 18).Value = 95.841656947549
$ws.Cells.Item(21, 19).Value = 0.007951090233125574
$ws.Cells.Item(21, 20).Value = 0.007951090233125572
$ws.Cells.Item(22, 7).Value = 11.40437333333333
$ws.Cells.Item(22, 8).Value = 34.21312
$ws.Cells.Item(22, 9).Value = 0.3307583551043392
$ws.Cells.Item(22, 10).Value = 0.3307583551043392
$ws.Cells.Item(22, 13).Value = 10.08846466666667
$ws.Cells.Item(22, 14).Value = 30.265394
$ws.Cells.Item(22, 15).Value = 0.2597171077778241
$ws.Cells.Item(22, 16).Value = 0.2597171077778241
$ws.Cells.Item(22, 17).Value = 115.0526174188089
$ws.Cells.Item(22, 18).Value = 1035.47355676928
$ws.Cells.Item(22, 19).Value = 0.08590360336104949
$ws.Cells.Item(22, 20).Value = 0.08590360336104948
$ws.Cells.Item(23, 7).Value = 11.40437333333333
$ws.Cells.Item(23, 8).Value = 34.21312
$ws.Cells.Item(23, 9).Value = 0.3307583551043392
$ws.Cells.Item(23, 10).Value = 0.3307583551043392
$ws.Cells.Item(23, 15).Value = 0.05876531725312483
$ws.Cells.Item(23, 16).Value = 0.05876531725312482
$ws.Cells.Item(23, 17).Value = 26.03256913365334
$ws.Cells.Item(23, 18).Value = 234.29312220288
$ws.Cells.Item(23, 19).Value = 0.01943711967182821
$ws.Cells.Item(23, 20).Value = 0.01943711967182821
$ws.Cells.Item(24, 7).Value = 11.40437333333333
$ws.Cells.Item(24, 8).Value = 34.21312
$ws.Cells.Item(24, 9).Value = 0.3307583551043392
$ws.Cells.Item(24, 10).Value = 0.3307583551043392
$ws.Cells.Item(24, 13).Value = 14.90894133333333
$ws.Cells.Item(24, 14).Value = 44.726824
$ws.Cells.Item(24, 15).Value = 0.3838153030278664
$ws.Cells.Item(24, 16).Value = 0.3838153030278664
$ws.Cells.Item(24, 17).Value = 170.0271329700978
$ws.Cells.Item(24, 18).Value = 1530.24419673088
$ws.Cells.Item(24, 19).Value = 0.1269501182933706
$ws.Cells.Item(24, 20).Value = 0.1269501182933706
$ws.Cells.Item(25, 7).Value = 11.40437333333333
$ws.Cells.Item(25, 8).Value = 34.21312
$ws.Cells.Item(25, 9).Value = 0.3307583551043392
$ws.Cells.Item(25, 10).Value = 0.3307583551043392
$ws.Cells.Item(25, 11).Value = 3.0
$ws.Cells.Item(25, 12).Value = 1.0
$ws.Cells.Item(25, 13).Value = 0.2836386666666667
$ws.Cells.Item(25, 14).Value = 0.850916
$ws.Cells.Item(25, 15).Value = 0.007301984652235982
$ws.Cells.Item(25, 16).Value = 0.007301984652235982
$ws.Cells.Item(25, 17).Value = 3.234721246435555
$ws.Cells.Item(25, 18).Value = 29.11249121792
$ws.Cells.Item(25, 19).Value = 0.002415192432570704
$ws.Cells.Item(25, 20).Value = 0.002415192432570704
$ws.Cells.Item(26, 7).Value = 11.40437333333333
$ws.Cells.Item(26, 8).Value = 34.21312
$ws.Cells.Item(26, 9).Value = 0.3307583551043392
$ws.Cells.Item(26, 10).Value = 0.3307583551043392
$ws.Cells.Item(26, 13).Value = 11.28032366666667
$ws.Cells.Item(26, 14).Value = 33.840971
$ws.Cells.Item(26, 15).Value = 0.2904002872889486
$ws.Cells.Item(26, 16).Value = 0.2904002872889486
$ws.Cells.Item(26, 17).Value = 128.6450224155022
$ws.Cells.Item(26, 18).Value = 1157.80520173952
$ws.Cells.Item(26, 19).Value = 0.09605232134552018
$ws.Cells.Item(26, 20).Value = 0.09605232134552016
